$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$newText = "This paper focuses on how COVID-19 information was communicated within and between different countries, reactions of governments to the pandemic, and attitudes and risk perceptions people had towards the virus. The major questions to answer are how digital communications influenced people’s interpretation of the news, what their responses were to the new laws and mandates, their beliefs and concerns about it versus other world issues, and the similarity and trends among the different countries."

# The paragraph is currently split across three runs. The first run
# ("This project will focus on how COVID-19 ", 40 characters) carries the
# formatting (solid black fill + Arial latin/ea/cs/sym) that should survive
# on the single merged run after the edit. Re-typing into that run's
# character range keeps its rPr intact (including the sym typeface, which
# isn't reachable through the Font object), then the leftover characters
# from the old runs 2 and 3 are deleted.
$firstRunLen = 40
$oldTotalLen = $tr.Length

$head = $tr.Characters(1, $firstRunLen)
$head.Text = $newText

$afterInsertLen = $shape.TextFrame.TextRange.Length
$newTextLen = $newText.Length
$remaining = $afterInsertLen - $newTextLen

if ($remaining -gt 0) {
    $tail = $shape.TextFrame.TextRange.Characters($newTextLen + 1, $remaining)
    $tail.Text = ""
}
